# Swap the two embedded themes' colour schemes.
#
# Before the edit:
#   ppt/theme/theme1.xml  = "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml  = "Integral"     (used by the Slide Master / the
#                                            deck's actual Design)
#
# After the edit the two files' contents are swapped, i.e. the deck's
# visible Design ("theme2.xml") now carries the stock "Office Theme"
# colour scheme while the colours that used to be on the Slide Master
# ("Integral") move to "theme1.xml".
#
# The only colour scheme reachable from the PowerPoint object model in
# this deck is the one backing the Slide Master / the slides themselves
# (ppt/theme/theme2.xml) - that's the "Integral" -> "Office Theme" half
# of the swap, so we apply it through ThemeColorScheme, which edits the
# 12 scheme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in
# the deck's active theme.

$p = $ppt.ActivePresentation

# Use slide 1's theme colour scheme - it resolves to the presentation's
# single active theme (ppt/theme/theme2.xml), which also backs
# $p.SlideMaster.ColorScheme / $p.SlideMaster.ThemeColorScheme.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Helper: pack R,G,B (0-255 each) into the COLORREF integer the RGB
# property expects (0x00BBGGRR).
function ColorRef([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours that used to live in
# ppt/theme/theme1.xml.
$tcs.Colors(1).RGB  = ColorRef 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = ColorRef 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = ColorRef 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = ColorRef 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = ColorRef 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = ColorRef 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = ColorRef 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = ColorRef 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = ColorRef 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = ColorRef 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = ColorRef 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = ColorRef 0x95 0x4F 0x72   # folHlink
